$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Present" count (value 3) for rows 14-21, column D ---
# These cells were previously blank; the diff adds a literal value of 3
# to each, while keeping their existing styles (s="25" / s="27") intact.
$ws.Range("D14").Value = 3
$ws.Range("D15").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("D17").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("D19").Value = 3
$ws.Range("D20").Value = 3
$ws.Range("D21").Value = 3

# --- Update the saved view state on the sheet ---
# Scroll the window so row 10 is the top-visible row, then move the
# active selection to I16 (previously H19).
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("I16").Select()
